# Strip the leading space from the text values of column B for the
# "category" rows (e.g. " Com banheiro ou sanitário de uso exclusivo dos
# moradores" -> "Com banheiro ou sanitário de uso exclusivo dos moradores").
#
# Affected rows: 7-33, 39-65, 71-97 (inclusive).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    @(7, 33),
    @(39, 65),
    @(71, 97)
)

foreach ($range in $ranges) {
    $startRow = $range[0]
    $endRow = $range[1]
    for ($row = $startRow; $row -le $endRow; $row++) {
        $cell = $ws.Cells.Item($row, 2)
        $value = $cell.Value2
        if ($null -ne $value -and $value -is [string] -and $value.StartsWith(" ")) {
            $cell.Value2 = $value.TrimStart(" ")
        }
    }
}
